# "fix the merge error" — a bad merge left several Public/Private/Save/View
# boolean flags on the Property sheet out of sync, and left the sheet
# scrolled/selected at the wrong cell. Restore the intended values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- Column F header is "View" (unchanged text, kept for clarity) ---
$ws.Range("F1").Value2 = "View"

# --- Rows 68-75: the "View" flag (column F) was dropped by the bad merge;
#     these property rows should all have View = TRUE ---
$ws.Range("F68").Value2 = $true
$ws.Range("F69").Value2 = $true
$ws.Range("F70").Value2 = $true
$ws.Range("F71").Value2 = $true
$ws.Range("F72").Value2 = $true
$ws.Range("F73").Value2 = $true
$ws.Range("F74").Value2 = $true
$ws.Range("F75").Value2 = $true

# --- Row 76 (GameID): Private/Save should be FALSE, View should be TRUE ---
$ws.Range("D76").Value2 = $false
$ws.Range("E76").Value2 = $false
$ws.Range("F76").Value2 = $true

# --- Row 77 (GateID): Private/Save should be FALSE, View should be TRUE ---
$ws.Range("D77").Value2 = $false
$ws.Range("E77").Value2 = $false
$ws.Range("F77").Value2 = $true

# --- Row 78 (GuildID): Public should be FALSE ---
$ws.Range("C78").Value2 = $false

# --- Restore the view/selection state of the sheet ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C78").Select() | Out-Null
